# 2nd batch JTown special picture links
#
# 1. Split the <dryingUrl> paragraph so the URL (plus the accidentally
#    over-selected "</dryingUrl" tail) becomes a real hyperlink, leaving the
#    opening "<dryingUrl>" tag and the trailing ">" as plain text runs.
# 2. Insert a blank paragraph followed by a new paragraph containing the
#    JTown "special picture" imgur link.
# 3. Register the (now-used) "Hyperlink" character style and the
#    "Unresolved Mention" character style in styles.xml, matching what Word
#    writes when it unhides/creates those styles.

$d = $word.ActiveDocument

# --- 1. Turn the drying-herbs URL into a hyperlink -------------------------

$dryingPara = $d.Paragraphs(3).Range
$searchRange = $d.Range($dryingPara.Start, $dryingPara.End)

$found = $searchRange.Find.Execute("https://www.thespruce.com/how-to-dry-and-store-herbs-1403397</dryingUrl", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Hyperlinks.Add($searchRange, "https://www.thespruce.com/how-to-dry-and-store-herbs-1403397") | Out-Null
}

# --- 2. Add the blank paragraph + the imgur picture-link paragraph ---------

$dryingPara = $d.Paragraphs(3).Range
$dryingPara.InsertParagraphAfter()

$blankPara = $d.Paragraphs(4).Range
$blankPara.InsertParagraphAfter()

$picturePara = $d.Paragraphs(5).Range
$picturePara.Text = "https://i.imgur.com/gpr7vky.jpg"

# --- 3. Make sure the character styles used/implied by the paste exist -----

$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = "DefaultParagraphFont"
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkStyle.Font.Color = 12673797
$hyperlinkStyle.Font.Underline = 1

$mentionStyle = $d.Styles.Add("Unresolved Mention", 2)
$mentionStyle.BaseStyle = "DefaultParagraphFont"
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.Color = 6053472
